$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 143, pushing the existing rows 143-177 down to 145-179
$ws.Range("143:144").Insert()

# New row 143
$ws.Cells.Item(143, 1).Value = 10
$ws.Cells.Item(143, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(143, 3).Value = "La Araucanía"
$ws.Cells.Item(143, 4).Value = 44964
$ws.Cells.Item(143, 5).Value = 9
$ws.Cells.Item(143, 6).Value = 100112031
$ws.Cells.Item(143, 7).Value = "Poroto verde"
$ws.Cells.Item(143, 8).Value = "Brío"
$ws.Cells.Item(143, 9).Value = "Primera"
$ws.Cells.Item(143, 10).Value = 200
$ws.Cells.Item(143, 11).Value = 1500
$ws.Cells.Item(143, 12).Value = 1500
$ws.Cells.Item(143, 13).Value = 1500
$ws.Cells.Item(143, 14).Value = "$/kilo"
$ws.Cells.Item(143, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(143, 16).Value = 1500
$ws.Cells.Item(143, 17).Value = 1
$ws.Cells.Item(143, 18).Value = "Hortaliza"

# New row 144
$ws.Cells.Item(144, 1).Value = 10
$ws.Cells.Item(144, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(144, 3).Value = "La Araucanía"
$ws.Cells.Item(144, 4).Value = 44964
$ws.Cells.Item(144, 5).Value = 9
$ws.Cells.Item(144, 6).Value = 100112031
$ws.Cells.Item(144, 7).Value = "Poroto verde"
$ws.Cells.Item(144, 8).Value = "Sin especificar"
$ws.Cells.Item(144, 9).Value = "Primera"
$ws.Cells.Item(144, 10).Value = 100
$ws.Cells.Item(144, 11).Value = 1500
$ws.Cells.Item(144, 12).Value = 1500
$ws.Cells.Item(144, 13).Value = 1500
$ws.Cells.Item(144, 14).Value = "$/kilo"
$ws.Cells.Item(144, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(144, 16).Value = 1500
$ws.Cells.Item(144, 17).Value = 1
$ws.Cells.Item(144, 18).Value = "Hortaliza"
